# Daily update at 8 AM UTC
# Adds the next day's row (row 38) to the "Wins Over Time" tracking sheet
# and rolls the date-only formatting down from the previous last row (37)
# to the new last row (38), matching the sheet's convention that only the
# bottom-most row's date cell uses the short "YYYY-MM-DD" format while all
# earlier rows use the full "YYYY-MM-DD HH:MM:SS" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 is no longer the last row, so give its date cell the same
# full-timestamp number format used by every other non-final row.
$ws.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 38: next day's data.
$ws.Range("A38").Value = 45622
$ws.Range("A38").NumberFormat = "YYYY-MM-DD"
$ws.Range("B38").Value = 97
$ws.Range("C38").Value = 77
$ws.Range("D38").Value = 90
